$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.842.65"
$ws.Range("E2").Value = "  +0.19%  "

# Row 3
$ws.Range("D3").Value = "3.812.33"
$ws.Range("E3").Value = "  +0.85%  "

# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").Value = "'605.06"
$ws.Range("E5").Value = "  +1.55%  "

# Row 6
$ws.Range("D6").Value = "'166.09"
$ws.Range("E6").Value = "  -0.70%  "

# Row 7
$ws.Range("E7").Value = "  +0.10%  "

# Row 8
$ws.Range("E8").Value = "  -0.08%  "

# Row 9
$ws.Range("E9").Value = "  +0.25%  "

# Row 10
$ws.Range("E10").Value = "  +0.88%  "

# Row 11
$ws.Range("E11").Value = "  +0.72%  "

# Row 13
$ws.Range("E13").Value = "  +0.03%  "

# Row 14
$ws.Range("D14").Value = "4.450.72"
$ws.Range("E14").Value = "  +0.84%  "

# Row 15
$ws.Range("D15").Value = "3.829.04"
$ws.Range("E15").Value = "  +1.28%  "

# Row 16
$ws.Range("D16").Value = "67.855.02"

# Row 17
$ws.Range("D17").Value = "'18.40"
$ws.Range("E17").Value = "  +0.07%  "

# Row 18
$ws.Range("E18").Value = "  +0.95%  "

# Row 19
$ws.Range("E19").Value = "  +1.77%  "

# Row 20
$ws.Range("D20").Value = "'463.94"
$ws.Range("E20").Value = "  +1.26%  "

# Row 22
$ws.Range("E22").Value = "  +0.94%  "

# Row 23
$ws.Range("E23").Value = "  -4.05%  "

# Row 24
$ws.Range("D24").Value = "'83.38"
$ws.Range("E24").Value = "  +0.15%  "

# Row 25
$ws.Range("E25").Value = "  +1.00%  "

# Row 26
$ws.Range("E26").Value = "  -0.71%  "

# Row 27
$ws.Range("D27").Value = "'10.03"
$ws.Range("E27").Value = "  -0.12%  "

# Row 28
$ws.Range("E28").Value = "  -0.11%  "

# Row 29
$ws.Range("D29").Value = "3.962.04"

# Row 30
$ws.Range("E30").Value = "  +0.99%  "

# Row 31
$ws.Range("E31").Value = "  +2.63%  "

# Row 32
$ws.Range("D32").Value = "'2.22"
$ws.Range("E32").Value = "  -1.29%  "

# Row 33
$ws.Range("D33").Value = "'29.58"
$ws.Range("E33").Value = "  -0.38%  "

# Row 34
$ws.Range("E34").Value = "  +0.12%  "

# Row 35
$ws.Range("E35").Value = "  -0.36%  "

# Row 36
$ws.Range("E36").Value = "  -0.04%  "

# Row 37
$ws.Range("E37").Value = "  +0.09%  "

# Row 38
$ws.Range("B38").Value = "Mantle"
$ws.Range("C38").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D38").Value = "'0.999"
$ws.Range("E38").Value = "  +0.34%  "

# Row 39
$ws.Range("B39").Value = "Filecoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D39").Value = "'5.82"
$ws.Range("E39").Value = "  +1.03%  "

# Row 40
$ws.Range("D40").Value = "'3.22"
$ws.Range("E40").Value = "  -4.39%  "

# Row 41
$ws.Range("E41").Value = "  -0.02%  "

# Row 43
$ws.Range("D43").Value = "'44.36"
$ws.Range("E43").Value = "  -2.95%  "

# Row 44
$ws.Range("D44").Value = "'47.70"
$ws.Range("E44").Value = "  -0.91%  "

# Row 45
$ws.Range("E45").Value = "  +0.61%  "

# Row 46
$ws.Range("D46").Value = "'28.05"
$ws.Range("E46").Value = "  +6.34%  "

# Row 47
$ws.Range("D47").Value = "'151.72"
$ws.Range("E47").Value = "  +1.62%  "

# Row 48
$ws.Range("E48").Value = "  +0.54%  "

# Row 49
$ws.Range("E49").Value = "  +11.50%  "

# Row 50
$ws.Range("E50").Value = "  +1.67%  "

# Row 51
$ws.Range("D51").Value = "'390.12"
$ws.Range("E51").Value = "  -0.66%  "
